$d = $word.ActiveDocument

# Locate the unique paragraph that begins the Executive Summary section:
# "Lorem ipsum dolor sit amet, consectetur adipiscing elit. Integer commodo..."
# (styled "FirstParagraph"). We match on a distinctive leading fragment of its
# text so we don't also match the later, differently-worded
# "New paragraph. Lorem ipsum..." paragraph found elsewhere in the document.
$targetIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t.StartsWith("Lorem ipsum dolor sit amet, consectetur adipiscing elit. Integer commodo gravida justo consectetur condimentum.")) {
        $targetIndex = $idx
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph"
}

$target = $d.Paragraphs($targetIndex)

# Split a brand-new paragraph in right before the target paragraph's text.
# Because we have not touched the target paragraph's style yet, the newly
# created paragraph inherits the existing "FirstParagraph" style, exactly as
# the diff requires (the original <w:pStyle w:val="FirstParagraph"/> stays
# attached to the new first paragraph).
$r = $target.Range.Duplicate
$r.Collapse(1)
$r.InsertBefore("Let us check if the documentation autoupdates with this added sentence.`r")

# The original Lorem-ipsum paragraph has now shifted one position later (it
# is no longer the first paragraph of the section); give it the explicit
# "BodyText" style called for by the diff.
$d.Paragraphs($targetIndex + 1).Style = "BodyText"
